$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new sheet "2022-Q1" right before the "总计" sheet.
#    Copying the existing "2021-Q4" sheet (identical B..H column
#    layout) gives us matching sheetPr/pageMargins/style plumbing.
# ------------------------------------------------------------------
$src   = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)

# The copy lands immediately before "总计"; re-fetch "总计" (a cached
# reference's .Index does not refresh in place) to locate the new sheet.
$total    = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($total.Index - 1)
$newSheet.Name = "2022-Q1"

# Restore the originally-active tab (copying makes the new sheet active).
$wb.Worksheets.Item(1).Activate()

# Header row (row 1)
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

function Set-TextValue($cell, $text) {
    # Force a text cell (matches the source data's inlineStr cells,
    # even for numeric-looking strings like "12.96") without leaving
    # a stray quote-prefix style behind.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2 — 161611 / 融通内需驱动混合
$newSheet.Cells.Item(2,1).Value = 0
Set-TextValue $newSheet.Cells.Item(2,2) "161611"
Set-TextValue $newSheet.Cells.Item(2,3) "融通内需驱动混合"
Set-TextValue $newSheet.Cells.Item(2,4) "12.96"
Set-TextValue $newSheet.Cells.Item(2,5) "65.58"
Set-TextValue $newSheet.Cells.Item(2,6) "2.28"
Set-TextValue $newSheet.Cells.Item(2,7) "0.2955"
$newSheet.Cells.Item(2,8).Value = 8

# Row 3 — 217024 / 招商安盈债券
$newSheet.Cells.Item(3,1).Value = 1
Set-TextValue $newSheet.Cells.Item(3,2) "217024"
Set-TextValue $newSheet.Cells.Item(3,3) "招商安盈债券"
Set-TextValue $newSheet.Cells.Item(3,4) "35.05"
Set-TextValue $newSheet.Cells.Item(3,5) "20.20"
Set-TextValue $newSheet.Cells.Item(3,6) "0.80"
Set-TextValue $newSheet.Cells.Item(3,7) "0.2804"
$newSheet.Cells.Item(3,8).Value = 7

# ------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing rows (and their running index in column A) down by one.
# ------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Rows.Item(2).ClearFormats()

# Re-apply column-A's numeric style (lost by ClearFormats) by copying
# the format only from the row right below it.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.58

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
